$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$combined = "('Nalathni Dragon', ['{2}{R}{R}', 'Creature " + [char]0x2014 + " Dragon', 'Flying; banding (Any creatures with banding, and up to one without, can attack in a band. Bands are blocked as a group. If any creatures with banding you control are blocking or being blocked by a creature, you divide that creature" + [char]0x2019 + "s combat damage, not its controller, among any of the creatures it" + [char]0x2019 + "s being blocked by or is blocking.)', '{R}: Nalathni Dragon gets +1/+0 until end of turn. If this ability has been activated four or more times this turn, sacrifice Nalathni Dragon at the beginning of the next end step.', '1/1'])"

$ws.Range("A2").Value = $combined

$ws.Range("A3:A7").EntireRow.Delete()
